$d = $word.ActiveDocument

# Create the new paragraph at the end of the document first (outside of
# track-changes) so the existing last paragraph is left completely
# untouched.
$newPara = $d.Paragraphs.Add()
$nr = $newPara.Range

# Turn on revision tracking only for the text insertion itself: this is
# the only reliable way to keep the newly typed text split across several
# <w:r> runs (matching how the text was originally authored/edited in
# separate passes) instead of Word's normal behaviour of silently
# collapsing adjacent runs that share identical formatting.
$d.TrackRevisions = $true

$nr.Collapse(0)
$nr.InsertAfter("Alcuni ticket presentano delle version non esistenti tra quelle elencate in Jira, questo capita")

$nr.Collapse(0)
$nr.InsertAfter(" soprattutto per i progetti ancora aperti. I")

$nr.Collapse(0)
$nr.InsertAfter("n tal caso si è scelto di settare la injected version come NULL")

$nr.Collapse(0)
$nr.InsertAfter(" e di scartarli in caso accadesse con la  opening e la fixed version")

$nr.Collapse(0)
$nr.InsertAfter(".")

$d.TrackRevisions = $false
[void]$d.AcceptAllRevisions()
